# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 13:27"

# Emiratos Arabes Unidos and Paises Bajos swap rank (rows 45/46)
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("A46").Value = "Paises Bajos"

# Nepal and Uzbekistan swap rank (rows 63/64)
$ws.Range("A63").Value = "Nepal"
$ws.Range("A64").Value = "Uzbekistan"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6291541
$ws.Range("C4").Value = 804
$ws.Range("D4").Value = 3547926
$ws.Range("E4").Value = 2553601
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 190014

# Row 15 - Iran
$ws.Range("B15").Value = 380746
$ws.Range("C15").Value = 1994
$ws.Range("D15").Value = 328595
$ws.Range("E15").Value = 30225
$ws.Range("G15").Value = 129
$ws.Range("H15").Value = 21926

# Row 23 - Alemania
$ws.Range("B23").Value = 247436
$ws.Range("C23").Value = 45
$ws.Range("E23").Value = 14943

# Row 44 - Bielorrusia
$ws.Range("B44").Value = 72302
$ws.Range("C44").Value = 161
$ws.Range("D44").Value = 71205
$ws.Range("E44").Value = 401
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 696

# Row 45 - Emiratos Arabes Unidos (updated values, now ranked above Paises Bajos)
$ws.Range("B45").Value = 72154
$ws.Range("C45").Value = 614
$ws.Range("D45").Value = 62668
$ws.Range("E45").Value = 9099
$ws.Range("H45").Value = 387

# Row 46 - Paises Bajos (values carried down, unchanged)
$ws.Range("B46").Value = 71863
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("H46").Value = 6235

# Row 62 - Suiza
$ws.Range("B62").Value = 43127
$ws.Range("C62").Value = 364
$ws.Range("E62").Value = 4615
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 2012

# Row 63 - Nepal (updated values, now ranked above Uzbekistan)
$ws.Range("B63").Value = 42877
$ws.Range("C63").Value = 1228
$ws.Range("D63").Value = 24207
$ws.Range("E63").Value = 18413
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 257

# Row 64 - Uzbekistan (values carried down, unchanged)
$ws.Range("B64").Value = 42540
$ws.Range("C64").Value = 103
$ws.Range("D64").Value = 39883
$ws.Range("E64").Value = 2328
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 329

# Row 84 - Libia
$ws.Range("B84").Value = 15773
$ws.Range("C84").Value = 617
$ws.Range("D84").Value = 1856
$ws.Range("E84").Value = 13663
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 254

# Row 85 - Madagascar
$ws.Range("B85").Value = 15106
$ws.Range("C85").Value = 83
$ws.Range("D85").Value = 14031
$ws.Range("E85").Value = 878
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 197

# Row 87 - Senegal
$ws.Range("B87").Value = 13826
$ws.Range("C87").Value = 83
$ws.Range("D87").Value = 9553
$ws.Range("E87").Value = 3986

# Row 130 - Eslovenia
$ws.Range("B130").Value = 3032
$ws.Range("C130").Value = 53
$ws.Range("D130").Value = 2402
$ws.Range("E130").Value = 496

# Row 146 - Malta
$ws.Range("B146").Value = 1965
$ws.Range("C146").Value = 34
$ws.Range("D146").Value = 1528
$ws.Range("E146").Value = 424

# Row 184 - Gibraltar
$ws.Range("B184").Value = 298
$ws.Range("C184").Value = 3
$ws.Range("D184").Value = 246
$ws.Range("E184").Value = 52
